# Swap the B:AC data (everything except the A "position" index column) between
# each of the following row pairs. The underlying match rows got re-sorted by
# a newer base-id (column B), but the positional index in column A was kept,
# so effectively the rest of each row's data (B and F..AC; C/D/E are identical
# between the two rows of a pair already) swaps places with its neighbour.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowPairs = @(
    @(62, 63),
    @(107, 108),
    @(125, 126),
    @(149, 150),
    @(172, 173),
    @(187, 188),
    @(191, 192),
    @(203, 204),
    @(225, 226),
    @(238, 239),
    @(243, 244)
)

foreach ($pair in $rowPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]

    $range1 = $ws.Range("B$r1`:AC$r1")
    $range2 = $ws.Range("B$r2`:AC$r2")

    $vals1 = $range1.Value2
    $vals2 = $range2.Value2

    $range1.Value2 = $vals2
    $range2.Value2 = $vals1
}
